$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash copies of the existing date/time cell formats in helper cells on row 1
# (row 1 is not touched by the row delete below, so the formats survive)
$ws.Range("B2").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C2").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats

# Remove the existing attendance entries (Yash, Ansh)
$ws.Rows("2:3").Delete()

# Re-apply the stashed date/time formats to new (blank) attendance rows
$ws.Range("E1").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("F1").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C13").PasteSpecial(-4122)

# Remove the helper cells used to stash the formats
$ws.Range("E1:F1").Clear()

# Move the active selection to the last entered cell
$null = $ws.Range("C13").Select()
